# Insert a new row 239 (shifts existing rows 239-282 down to 240-283)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(239).Insert()

# Populate the newly inserted row 239 with the new weekly price record
$ws.Range("A239").Value = 11
$ws.Range("B239").Value = "Vega Monumental Concepción"
$ws.Range("C239").Value = "Bíobío"
$ws.Range("D239").Value = 45218
$ws.Range("E239").Value = 8
$ws.Range("F239").Value = "Fruta"
$ws.Range("G239").Value = 100102
$ws.Range("H239").Value = "Cítricos"
$ws.Range("I239").Value = 100102004
$ws.Range("J239").Value = "Mandarina"
$ws.Range("K239").Value = "Murcott"
$ws.Range("L239").Value = "Primera"
$ws.Range("M239").Value = 180
$ws.Range("N239").Value = 8000
$ws.Range("O239").Value = 9000
$ws.Range("P239").Value = 8444
$ws.Range("Q239").Value = "$/bandeja 18 kilos"
$ws.Range("R239").Value = "Región de O'Higgins"
$ws.Range("S239").Value = 469
$ws.Range("T239").Value = 18
